# temp.m - progressing on constraint definition; no longer using cvx
$wb = $excel.ActiveWorkbook

# --- Remove the unused "Sheet1" worksheet ---
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Delete()

# --- Add new columns (kA, backgroundConc) to the "m" sheet ---
$ws = $wb.Worksheets.Item("m")

$ws.Range("E1").Value = "kA"
$ws.Range("F1").Value = "backgroundConc"

$ws.Range("E2").Value = 0.14
$ws.Range("F2").Value = 3

$ws.Range("E3").Value = 0.074
$ws.Range("F3").Value = 1

$ws.Range("E4").Value = 0.074
$ws.Range("F4").Value = 1

# left-align the new numeric cells like the rest of the sheet
$ws.Range("E2:F4").HorizontalAlignment = -4131

# highlight E4 with a yellow fill
$ws.Range("E4").Interior.Color = 65535

$ws.Range("G2").Select()
